$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (id=1, "Edwards-Thomas Incubator") entirely; subsequent rows
# shift up and the sheet's used range shrinks from A1:H31 to A1:H30.
$ws.Rows.Item(2).Delete()
